# B1--and-B2-PowerPoint.pptx : Fri, Jun 26, 2020 10:05:02 PM
#
# 1) The table on slide 5 gets a new table style (tableStyleId change).
# 2) The deck's design theme is switched from the "Integral" / "Red Violet"
#    colour scheme to the standard "Office Theme" / "Office" colour scheme.

$p = $ppt.ActivePresentation

# --- 1) Re-style the financial-documents table on slide 5 -------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{73AF0AB2-5345-4D0C-85A4-359BBDBCB833}")

# --- 2) Swap the theme colours back to the default Office palette -----------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeRGB {
    param(
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $oleColor = $r + ($g * 256) + ($b * 65536)
    $themeColors.Colors($Index).RGB = $oleColor
}

Set-ThemeRGB 1  "000000"   # dk1
Set-ThemeRGB 2  "FFFFFF"   # lt1
Set-ThemeRGB 3  "44546A"   # dk2
Set-ThemeRGB 4  "E7E6E6"   # lt2
Set-ThemeRGB 5  "5B9BD5"   # accent1
Set-ThemeRGB 6  "ED7D31"   # accent2
Set-ThemeRGB 7  "A5A5A5"   # accent3
Set-ThemeRGB 8  "FFC000"   # accent4
Set-ThemeRGB 9  "4472C4"   # accent5
Set-ThemeRGB 10 "70AD47"   # accent6
Set-ThemeRGB 11 "0563C1"   # hlink
Set-ThemeRGB 12 "954F72"   # folHlink
